$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated s_val data (regenerated to filter save games)
$ws.Range("B2").Value = 1.505614041169197
$ws.Range("C2").Value = 1.65323645889881
$ws.Range("D2").Value = 0.7127328510149897
$ws.Range("E2").Value = 0.4998867070740569
$ws.Range("G2").Value = 4.371470058157054

$ws.Range("B3").Value = 1.505614041169197
$ws.Range("C3").Value = 1.65323645889881
$ws.Range("D3").Value = 0.7127328510149897
$ws.Range("E3").Value = 0.4998867070740569
$ws.Range("G3").Value = 4.371470058157054

$ws.Range("B4").Value = 3.182878228561681
$ws.Range("C4").Value = 1.65323645889881
$ws.Range("D4").Value = 0.7127328510149897
$ws.Range("E4").Value = 0.4998867070740569
$ws.Range("G4").Value = 6.048734245549538

$ws.Range("B5").Value = 1.505614041169197
$ws.Range("C5").Value = 1.65323645889881
$ws.Range("D5").Value = 3.082599426703578
$ws.Range("E5").Value = 0.4998867070740569
$ws.Range("G5").Value = 6.741336633845642

$ws.Range("B6").Value = 0.3464964993005633
$ws.Range("C6").Value = 0.004309184025731883
$ws.Range("D6").Value = 3.082599426703578
$ws.Range("E6").Value = 0.4998867070740569
$ws.Range("G6").Value = 3.933291817103931
